# Generate Report for Handback
#
# - Status column (C) on the zh-cn / de-de sheets moves from "Ready for
#   handoff" to "Handed back: in sync with en-US" for both data rows.
# - The "Latest Handback DateTime" column (H) is stamped with the actual
#   handback timestamp (different per locale).
# - The "Latest Target File" (F) and "Latest Handback File" (G) columns are
#   now populated with hyperlinks to the source markdown file and the
#   handed-back xlf file, mirroring the existing hyperlink look (underlined,
#   cornflower-blue) already used by columns A and D.

$wb = $excel.ActiveWorkbook

# Cornflower blue (#6495ED) expressed the way the COM Color setter expects
# (it reads the integer as 0x00BBGGRR, the same packing VBA's RGB() uses).
$hyperlinkColor = 15570276

function Style-AsHyperlink($range) {
    $range.Font.Name = "Calibri"
    $range.Font.Size = 11
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsZh.Range("H2").Value = "2016-03-23 18:35:53"
$wsZh.Range("H3").Value = "2016-03-23 18:35:53"

$zhMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/6b693820c2c6973ef23b0f17b07748078663a309/e2e/a.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f21b87391755c57f9520d9fb0822a90a1b54554e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhXlfUrl, [Type]::Missing, [Type]::Missing, $zhXlfName)
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhXlfUrl, [Type]::Missing, [Type]::Missing, $zhXlfName)

Style-AsHyperlink $wsZh.Range("F2")
Style-AsHyperlink $wsZh.Range("G2")
Style-AsHyperlink $wsZh.Range("F3")
Style-AsHyperlink $wsZh.Range("G3")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe.Range("H2").Value = "2016-03-23 18:36:02"
$wsDe.Range("H3").Value = "2016-03-23 18:36:02"

$deMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/6b693820c2c6973ef23b0f17b07748078663a309/e2e/a.md"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/730e6c818a20611bcd19a634b3294253b1208093/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$deXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deXlfUrl, [Type]::Missing, [Type]::Missing, $deXlfName)
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deXlfUrl, [Type]::Missing, [Type]::Missing, $deXlfName)

Style-AsHyperlink $wsDe.Range("F2")
Style-AsHyperlink $wsDe.Range("G2")
Style-AsHyperlink $wsDe.Range("F3")
Style-AsHyperlink $wsDe.Range("G3")
